$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.2217898832684825
$ws.Range("C2").Value = 0.5136186770428015
$ws.Range("J2").Value = 0.003891050583657588
$ws.Range("O2").Value = 0.003891050583657588
$ws.Range("P2").Value = 0.1556420233463035
$ws.Range("S2").Value = 0.1011673151750973
$ws.Range("B3").Value = 0.0072992700729927
$ws.Range("C3").Value = 0.0291970802919708
$ws.Range("P3").Value = 0.7299270072992701
$ws.Range("S3").Value = 0.2335766423357664
$ws.Range("J4").Value = 0.05714285714285714
$ws.Range("P4").Value = 0.7142857142857143
$ws.Range("S4").Value = 0.2285714285714286
$ws.Range("B6").Value = 0.07389162561576355
$ws.Range("D6").Value = 0.009852216748768473
$ws.Range("F6").Value = 0.04433497536945813
$ws.Range("J6").Value = 0.3054187192118227
$ws.Range("O6").Value = 0.02463054187192118
$ws.Range("Q6").Value = 0.1773399014778325
$ws.Range("R6").Value = 0.1379310344827586
$ws.Range("S6").Value = 0.2266009852216749
$ws.Range("B7").Value = 0.07228915662650602
$ws.Range("D7").Value = 0.006024096385542169
$ws.Range("F7").Value = 0.03614457831325301
$ws.Range("J7").Value = 0.1024096385542169
$ws.Range("O7").Value = 0.01807228915662651
$ws.Range("Q7").Value = 0.1927710843373494
$ws.Range("R7").Value = 0.1144578313253012
$ws.Range("S7").Value = 0.4578313253012048
$ws.Range("B8").Value = 0.06506024096385542
$ws.Range("D8").Value = 0.02409638554216868
$ws.Range("F8").Value = 0.06024096385542169
$ws.Range("J8").Value = 0.1180722891566265
$ws.Range("O8").Value = 0.01204819277108434
$ws.Range("Q8").Value = 0.2289156626506024
$ws.Range("R8").Value = 0.1132530120481928
$ws.Range("S8").Value = 0.3783132530120482
$ws.Range("B9").Value = 0.1038961038961039
$ws.Range("D9").Value = 0.0303030303030303
$ws.Range("F9").Value = 0.05627705627705628
$ws.Range("J9").Value = 0.07792207792207792
$ws.Range("O9").Value = 0.02164502164502164
$ws.Range("Q9").Value = 0.2164502164502164
$ws.Range("R9").Value = 0.09956709956709957
$ws.Range("S9").Value = 0.3939393939393939
$ws.Range("B10").Value = 0.09001512859304085
$ws.Range("D10").Value = 0.01285930408472012
$ws.Range("F10").Value = 0.05824508320726172
$ws.Range("J10").Value = 0.1043872919818457
$ws.Range("O10").Value = 0.01210287443267776
$ws.Range("Q10").Value = 0.2329803328290469
$ws.Range("R10").Value = 0.1021180030257186
$ws.Range("S10").Value = 0.3872919818456884
$ws.Range("G11").Value = 0.1333333333333333
$ws.Range("J11").Value = 0.07916666666666666
$ws.Range("K11").Value = 0.1916666666666667
$ws.Range("L11").Value = 0.5875
$ws.Range("S11").Value = 0.008333333333333333
$ws.Range("G12").Value = 0.8082191780821918
$ws.Range("J12").Value = 0.1438356164383562
$ws.Range("K12").Value = 0.00684931506849315
$ws.Range("L12").Value = 0.0410958904109589
$ws.Range("F15").Value = 0.03896103896103896
$ws.Range("H15").Value = 0.1168831168831169
$ws.Range("I15").Value = 0.06493506493506493
$ws.Range("J15").Value = 0.4155844155844156
$ws.Range("K15").Value = 0.03463203463203463
$ws.Range("M15").Value = 0.008658008658008658
$ws.Range("O15").Value = 0.06493506493506493
$ws.Range("S15").Value = 0.2554112554112554
$ws.Range("F16").Value = 0.00625
$ws.Range("H16").Value = 0.14375
$ws.Range("I16").Value = 0.1125
$ws.Range("J16").Value = 0.48125
$ws.Range("K16").Value = 0.0625
$ws.Range("M16").Value = 0.0125
$ws.Range("O16").Value = 0.09375
$ws.Range("S16").Value = 0.08749999999999999
$ws.Range("F17").Value = 0.01532567049808429
$ws.Range("H17").Value = 0.1800766283524904
$ws.Range("I17").Value = 0.09386973180076628
$ws.Range("J17").Value = 0.4348659003831418
$ws.Range("K17").Value = 0.06896551724137931
$ws.Range("M17").Value = 0.01532567049808429
$ws.Range("O17").Value = 0.07854406130268199
$ws.Range("S17").Value = 0.1130268199233716
$ws.Range("F18").Value = 0.012
$ws.Range("H18").Value = 0.128
$ws.Range("I18").Value = 0.144
$ws.Range("J18").Value = 0.4
$ws.Range("K18").Value = 0.08799999999999999
$ws.Range("M18").Value = 0.004
$ws.Range("O18").Value = 0.08
$ws.Range("S18").Value = 0.144
$ws.Range("F19").Value = 0.01930812550281577
$ws.Range("H19").Value = 0.1898632341110217
$ws.Range("I19").Value = 0.09090909090909091
$ws.Range("J19").Value = 0.4070796460176991
$ws.Range("K19").Value = 0.09332260659694289
$ws.Range("M19").Value = 0.01367658889782784
$ws.Range("N19").Value = 0.001609010458567981
$ws.Range("O19").Value = 0.06355591311343524
$ws.Range("S19").Value = 0.1206757843925985
